# Depersonalize the data: replace the identifying index in column A with a
# simple sequential counter (0, 1, 2, ...) and update column B (diagnosis)
# to the corresponding values after the re-indexing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bValues = @(1,1,2,2,2,2,1,1,0,0,1,1,3,3,1,1,2,2,1,1,1,3,3,2,2,3,3,2,2,1,1,1,1,2,2,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,2,2,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,2,2,2,2,2,1,1,1,1,1,1,1,1,2,2,1,1,1,1,2,2,1,1,2,2,2,2,2,2,1,1,2,2,1,1,2,1,1,1,1,1,1,1,1,1,1,2,2,2,2,2,2,2,1,2,2,1,2,2,1,1,1,1,1,1,1,1,1)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
